# Update the date line and the division problems in the table.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-15 Thursday" "2025-05-16 Friday"

Replace-Text "766÷8=" "787÷3="
Replace-Text "589÷4=" "356÷6="
Replace-Text "697÷3=" "342÷2="
Replace-Text "180÷2=" "934÷6="
Replace-Text "232÷6=" "364÷2="

Replace-Text "481÷7=" "392÷4="
Replace-Text "831÷4=" "783÷6="
Replace-Text "133÷3=" "295÷2="
Replace-Text "175÷7=" "691÷8="
Replace-Text "214÷2=" "370÷3="

Replace-Text "133÷9=" "910÷4="
Replace-Text "229÷8=" "328÷5="
Replace-Text "578÷7=" "444÷6="
Replace-Text "255÷6=" "872÷3="
Replace-Text "357÷3=" "873÷8="

Replace-Text "181÷7=" "139÷2="
Replace-Text "957÷3=" "569÷7="
Replace-Text "479÷7=" "758÷8="
Replace-Text "113÷5=" "183÷3="
Replace-Text "522÷8=" "795÷8="

Replace-Text "436÷5=" "789÷8="
Replace-Text "720÷9=" "806÷9="
Replace-Text "197÷9=" "478÷2="
Replace-Text "695÷3=" "727÷2="
Replace-Text "303÷5=" "474÷2="
